# Weekly CompStat update (33rd Precinct) -- "New crime data collected"
# Applies the numeric refresh for the week of 3/31/2025-4/6/2025 (Volume 32, Number 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Banner text: issue number + reporting week dates ---
$ws.Range("A8").Value = "Volume 32   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  3/31/2025  Through  4/6/2025"

# --- Refreshed weekly/28-day/YTD/2-year crime counts + % changes ---
$ws.Range("M15").Value = 33.333333333333
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = -61.111111111111
$ws.Range("I16").Value = 23
$ws.Range("J16").Value = 44
$ws.Range("K16").Value = -47.727272727272
$ws.Range("L16").Value = -41.025641025641
$ws.Range("M16").Value = -57.407407407407
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 26.666666666666
$ws.Range("I17").Value = 51
$ws.Range("J17").Value = 65
$ws.Range("K17").Value = -21.538461538461
$ws.Range("L17").Value = -19.047619047619
$ws.Range("M17").Value = 8.510638297872
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 42.857142857142
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 33
$ws.Range("K18").Value = -27.272727272727
$ws.Range("L18").Value = -47.826086956521
$ws.Range("M18").Value = -7.692307692307
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -50
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = -51.515151515151
$ws.Range("I19").Value = 79
$ws.Range("J19").Value = 96
$ws.Range("K19").Value = -17.708333333333
$ws.Range("L19").Value = 3.947368421052
$ws.Range("M19").Value = 31.666666666666
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -50
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = 12
$ws.Range("L20").Value = -24.324324324324
$ws.Range("M20").Value = 40
$ws.Range("C21").Value = 17
$ws.Range("E21").Value = -10.526315789473
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = -17.5
$ws.Range("I21").Value = 209
$ws.Range("J21").Value = 266
$ws.Range("K21").Value = -21.428571428571
$ws.Range("L21").Value = -21.132075471698
$ws.Range("M21").Value = -1.415094339622
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("L23").Value = 33.333333333333
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 10.526315789473
$ws.Range("F24").Value = 62
$ws.Range("G24").Value = 70
$ws.Range("H24").Value = -11.428571428571
$ws.Range("I24").Value = 236
$ws.Range("J24").Value = 208
$ws.Range("K24").Value = 13.461538461538
$ws.Range("L24").Value = -3.673469387755
$ws.Range("M24").Value = 98.319327731092
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -61.538461538461
$ws.Range("I25").Value = 64
$ws.Range("J25").Value = 41
$ws.Range("K25").Value = 56.097560975609
$ws.Range("L25").Value = -4.477611940298
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 133.333333333333
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 12
$ws.Range("I26").Value = 86
$ws.Range("J26").Value = 85
$ws.Range("K26").Value = 1.176470588235
$ws.Range("L26").Value = -22.522522522522
$ws.Range("M26").Value = -14.851485148514
$ws.Range("L27").Value = -20
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 250
$ws.Range("I28").Value = 19
$ws.Range("K28").Value = 72.727272727272
$ws.Range("L28").Value = 111.111111111111

# --- Cells toggling between a numeric count and the "0"/"***.*" (N/A) text markers ---
# Housing (row 23): now has zero complaints both years -> "0" / "0" / "***.*"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E23").PasteSpecial(-4122)

# Other Sex Crimes (row 28): prior-year column now "0"/"N.A" instead of current-year
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("F23").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 2

